$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.081.89'
$ws.Range('E2').Value = '  +1.73%  '
$ws.Range('D3').Value = '3.775.10'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '628.04'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +4.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.44'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +1.37%  '
$ws.Range('D7').Value = '3.774.43'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +1.36%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  +2.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.78'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -0.78%  '
$ws.Range('E13').Value = '  -0.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.41'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.15%  '
$ws.Range('D15').Value = '4.408.82'
$ws.Range('E15').Value = '  -0.37%  '
$ws.Range('D16').Value = '3.785.98'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '69.091.26'
$ws.Range('E17').Value = '  +1.74%  '
$ws.Range('E18').Value = '  -3.10%  '
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.05'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '467.67'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.56'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.706'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.99'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('E25').Value = '  +0.92%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.02'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +1.36%  '
$ws.Range('E27').Value = '  +3.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.03'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.30%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '3.923.87'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('E31').Value = '  +3.01%  '
$ws.Range('E32').Value = '  +2.45%  '
$ws.Range('E33').Value = '  -0.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.76'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -0.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.172'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +18.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = '3.724.87'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.95'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.967'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '154.42'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '43.23'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('B47').Value = 'TheGraph'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.295'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +0.69%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.76'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -0.74%  '
$ws.Range('E49').Value = '  +3.91%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.38'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +1.30%  '
$ws.Range('E51').Value = '  -1.09%  '
